# Auto-generated Excel COM-interop script
# Applies the numeric "want-to-go" (F) / "min price" (G) counter bumps
# and one refreshed cover-image URL (I45), mirrored across the
# "展览" / "演出" / "全部类型" sheets exactly as the source diff shows.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

$wsExpo.Range("F4").Value = 63
$wsExpo.Range("F5").Value = 65
$wsExpo.Range("F6").Value = 820
$wsExpo.Range("F7").Value = 398
$wsExpo.Range("F8").Value = 4663
$wsExpo.Range("F9").Value = 4663
$wsExpo.Range("F12").Value = 153
$wsExpo.Range("F15").Value = 109
$wsExpo.Range("F16").Value = 7380
$wsExpo.Range("F21").Value = 508
$wsExpo.Range("F22").Value = 1342
$wsExpo.Range("F24").Value = 6281
$wsExpo.Range("F25").Value = 1733
$wsExpo.Range("G25").Value = 70
$wsExpo.Range("F26").Value = 18
$wsExpo.Range("F28").Value = 6146
$wsExpo.Range("F29").Value = 140
$wsExpo.Range("F31").Value = 113
$wsExpo.Range("F34").Value = 6363
$wsExpo.Range("F35").Value = 22
$wsExpo.Range("F39").Value = 18
$wsExpo.Range("F41").Value = 2447
$wsExpo.Range("F43").Value = 55
$wsExpo.Range("I45").Value = "//i0.hdslb.com/bfs/openplatform/202403/DHBY1mGz1711355939240.jpeg"
$wsExpo.Range("F46").Value = 421
$wsExpo.Range("F47").Value = 2130
$wsExpo.Range("F48").Value = 40
$wsShow.Range("F3").Value = 227
$wsShow.Range("F6").Value = 119
$wsShow.Range("F10").Value = 6
$wsShow.Range("F13").Value = 145
$wsAll.Range("F4").Value = 63
$wsAll.Range("F5").Value = 227
$wsAll.Range("F6").Value = 65
$wsAll.Range("F8").Value = 398
$wsAll.Range("F9").Value = 4663
$wsAll.Range("F10").Value = 4663
$wsAll.Range("F13").Value = 153
$wsAll.Range("F16").Value = 109
$wsAll.Range("F17").Value = 7380
$wsAll.Range("F20").Value = 508
$wsAll.Range("F21").Value = 1342
$wsAll.Range("F22").Value = 119
$wsAll.Range("F23").Value = 6281
$wsAll.Range("F24").Value = 1733
$wsAll.Range("G24").Value = 70
$wsAll.Range("F29").Value = 6146
$wsAll.Range("F30").Value = 140
$wsAll.Range("F31").Value = 6
$wsAll.Range("F33").Value = 113
$wsAll.Range("F36").Value = 6363
$wsAll.Range("F37").Value = 22
$wsAll.Range("F40").Value = 18
$wsAll.Range("F42").Value = 2447
$wsAll.Range("I45").Value = "//i0.hdslb.com/bfs/openplatform/202403/DHBY1mGz1711355939240.jpeg"
$wsAll.Range("F46").Value = 421
$wsAll.Range("F47").Value = 145
$wsAll.Range("F48").Value = 2130
$wsAll.Range("F49").Value = 40

Write-Output "Applied 59 cell updates across 展览 / 演出 / 全部类型."
